# Generate Report for Handback
# Updates the localization-status workbook after a handback: the overall
# status moves from "Ready for handoff" to "Handed back: in sync with en-US",
# and the per-language sheets (zh-cn / de-de) record the returned target
# file, the handback xliff file, and the handback timestamp.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$mdFileName = "a5cfe79f-3c5f-4df4-8006-b4782083b8ec.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a02f2c01a76ca446607a0ef475c3500a7439b9bc/e2e/a5cfe79f-3c5f-4df4-8006-b4782083b8ec.md"

# --- Overview sheet: update the zh-cn / de-de status columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E:F").AutoFit()

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("I2").Value = $mdFileName
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl, "", "", $mdFileName)
$wsZh.Range("J2").Value = "a5cfe79f-3c5f-4df4-8006-b4782083b8ec.a0fe7ee353b6f86ce7cb473353175e76f8d7bec5.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-07 09:31:53"
$wsZh.Range("C:C").AutoFit()
$wsZh.Range("I:J").AutoFit()

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("I2").Value = $mdFileName
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl, "", "", $mdFileName)
$wsDe.Range("J2").Value = "a5cfe79f-3c5f-4df4-8006-b4782083b8ec.a0fe7ee353b6f86ce7cb473353175e76f8d7bec5.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-07 09:32:01"
$wsDe.Range("C:C").AutoFit()
$wsDe.Range("I:J").AutoFit()
